$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: updated bird/cage info. Introduce the brand-new strings first
# (Australian Gouldian, Coastal Cities) so the shared-string table ends up
# ordered the same way the original author's workbook has it.
$ws.Range("B8").Value = "Australian Gouldian"
$ws.Range("C8").Value = "Coastal Cities"

# Row 5: SubSpec changed from "50A" to "13A"
$ws.Range("D5").Value = "13A"

# Row 7: SubSpec stays textually "15A" (shared string index shifts upstream, no
# visible change needed, but set explicitly for safety / parity)
$ws.Range("D7").Value = "15A"

$ws.Range("A8").Value = 101
$ws.Range("D8").Value = "12A"
$ws.Range("G8").Value = 111
$ws.Range("H8").ClearFormats()
$ws.Range("H8").Value = "15/05/2023"
$ws.Range("I8").Value = "Black"
$ws.Range("K8").Value = "Green"

# Row 9: removed entirely (data now ends at row 8)
$ws.Range("A9:K9").Delete()

# Restore cursor/selection position as left by the author
[void]$ws.Range("L6").Select()
